$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Move the "_GoBack" bookmark from right after "Instructions" in the title
#    (where it currently sits) down into the paragraph that, after step 2,
#    will be the now-empty paragraph following "cd <PROJECT PATH>\mdcs".
#    Re-adding a bookmark with an existing name relocates it, so we do not
#    need to explicitly delete the old one first.
# ---------------------------------------------------------------------------

# ---------------------------------------------------------------------------
# 2) Remove the "mkdir data" / "cd data" / "mkdir ts" / "mkdir db" paragraphs
#    entirely (their content and their paragraph marks), leaving the empty
#    paragraph that used to follow them intact.
# ---------------------------------------------------------------------------
$findStart = $d.Content
$null = $findStart.Find.Execute("mkdir data", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$startPara = $findStart.Paragraphs(1)
$startPos = $startPara.Range.Start

$findEnd = $d.Content
$null = $findEnd.Find.Execute("mkdir db", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$endPara = $findEnd.Paragraphs(1)
$endPos = $endPara.Range.End

$deleteRange = $d.Range($startPos, $endPos)
$deleteRange.Delete()

# After the deletion, $startPos now sits at the start of the paragraph that
# used to immediately follow "mkdir db" (an empty paragraph) -- that is
# where the bookmark belongs now.
$targetPara = $d.Range($startPos, $startPos).Paragraphs(1)
$d.Bookmarks.Add("_GoBack", $targetPara.Range)

# ---------------------------------------------------------------------------
# 3) Merge the two runs that hold the Jena "Option 2" command line (removing
#    the lastRenderedPageBreak between them) by re-writing the text that
#    spans the run boundary with itself, which causes Word to coalesce the
#    runs into a single one.
# ---------------------------------------------------------------------------
$mergeRange = $d.Content
$null = $mergeRange.Find.Execute("127.0.0.1:5556", $true, $false, $false, $false, $false, $true, 1, $false, "127.0.0.1:5556", 2)
